# Update column G ("K") values on the active sheet with newly computed
# values (regenerated using K instead of Strike#, per commit message).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 3
    3  = 6
    4  = 10
    5  = 5
    6  = 0
    7  = 2
    8  = 5
    9  = 2
    10 = 3
    11 = 6
    12 = 2
    13 = 3
    14 = 5
    15 = 5
    16 = 5
    17 = 6
    18 = 7
    19 = 6
    20 = 6
    21 = 6
    22 = 5
    23 = 2
    24 = 4
    25 = 6
    26 = 3
    27 = 4
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
